# Apply updated "dSF" (column F) values that were repulled from source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 7
    11 = 1
    16 = 0
    17 = -3
    19 = -1
    21 = -3
    28 = -1
    31 = -12
    35 = 7
    37 = -4
    38 = -2
    46 = 0
    47 = -1
    50 = -6
    52 = 7
    53 = 6
    54 = -5
    55 = -1
    56 = -2
    57 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
